# chore: update Sheets via scheduled runner
# Refresh of market-price-derived columns (currentAveragePrice*, LevePrice*,
# LeveProfit*) across the per-class Leve tables, rows keyed by Leve Item ID.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 445.85715
$ws.Range("I2").Value = 355.25
$ws.Range("J2").Value = 566.6667
$ws.Range("K2").Value = 355.25
$ws.Range("L2").Value = 566.6667
$ws.Range("M2").Value = -242.25
$ws.Range("N2").Value = -792.6667

$ws.Range("H33").Value = 1506.0952
$ws.Range("I33").Value = 564.375
$ws.Range("J33").Value = 4519.6
$ws.Range("K33").Value = 564.375
$ws.Range("L33").Value = 4519.6
$ws.Range("M33").Value = -335.375
$ws.Range("N33").Value = -4977.6

$ws.Range("H39").Value = 992470.4
$ws.Range("I39").Value = 1323028
$ws.Range("J39").Value = 797.5
$ws.Range("K39").Value = 3969084
$ws.Range("L39").Value = 2392.5
$ws.Range("M39").Value = -3968788
$ws.Range("N39").Value = -2984.5

$ws.Range("H40").Value = 2595.2083
$ws.Range("I40").Value = 2730
$ws.Range("J40").Value = 2435.9092
$ws.Range("K40").Value = 2730
$ws.Range("L40").Value = 2435.9092
$ws.Range("M40").Value = -2555
$ws.Range("N40").Value = -2785.9092

$ws.Range("H75").Value = 23314
$ws.Range("J75").Value = 23314
$ws.Range("L75").Value = 23314
$ws.Range("N75").Value = -25186

$ws.Range("H78").Value = 23314
$ws.Range("J78").Value = 23314
$ws.Range("L78").Value = 69942
$ws.Range("N78").Value = -79302

$ws.Range("H132").Value = 5106672.5
$ws.Range("I132").Value = 5686907.5
$ws.Range("J132").Value = 604.6
$ws.Range("K132").Value = 17060722.5
$ws.Range("L132").Value = 1813.8
$ws.Range("M132").Value = -17058192.5
$ws.Range("N132").Value = -6873.8

$ws.Range("H137").Value = 1897.8334
$ws.Range("I137").Value = 1498.7858
$ws.Range("J137").Value = 3294.5
$ws.Range("K137").Value = 4496.357400000001
$ws.Range("L137").Value = 9883.5
$ws.Range("M137").Value = -1946.357400000001
$ws.Range("N137").Value = -14983.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2681.7407
$ws.Range("I61").Value = 2673.2856
$ws.Range("J61").Value = 2684.7
$ws.Range("K61").Value = 2673.2856
$ws.Range("L61").Value = 2684.7
$ws.Range("M61").Value = -2461.2856
$ws.Range("N61").Value = -3108.7

$ws.Range("H69").Value = 34184.855
$ws.Range("J69").Value = 34184.855
$ws.Range("L69").Value = 34184.855
$ws.Range("N69").Value = -35682.855

$ws.Range("H72").Value = 34184.855
$ws.Range("J72").Value = 34184.855
$ws.Range("L72").Value = 102554.565
$ws.Range("N72").Value = -110042.565

$ws.Range("H136").Value = 2681.7407
$ws.Range("I136").Value = 2673.2856
$ws.Range("J136").Value = 2684.7
$ws.Range("K136").Value = 8019.8568
$ws.Range("L136").Value = 8054.099999999999
$ws.Range("M136").Value = -5469.8568
$ws.Range("N136").Value = -13154.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 841
$ws.Range("I94").Value = 761.5
$ws.Range("K94").Value = 761.5
$ws.Range("M94").Value = -310.5

$ws.Range("H134").Value = 2467.761
$ws.Range("I134").Value = 2503.3076
$ws.Range("J134").Value = 2269.7144
$ws.Range("K134").Value = 7509.9228
$ws.Range("L134").Value = 6809.1432
$ws.Range("M134").Value = -4974.9228
$ws.Range("N134").Value = -11879.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 11820.25
$ws.Range("I58").Value = 2151.25
$ws.Range("J58").Value = 21489.25
$ws.Range("K58").Value = 2151.25
$ws.Range("L58").Value = 21489.25
$ws.Range("M58").Value = -1948.25
$ws.Range("N58").Value = -21895.25

$ws.Range("H132").Value = 3462.6511
$ws.Range("I132").Value = 3328.8823
$ws.Range("J132").Value = 3968
$ws.Range("K132").Value = 9986.6469
$ws.Range("L132").Value = 11904
$ws.Range("M132").Value = -7456.6469
$ws.Range("N132").Value = -16964

$ws.Range("H134").Value = 1180.6061
$ws.Range("I134").Value = 1180.6061
$ws.Range("K134").Value = 3541.8183
$ws.Range("M134").Value = -1006.8183

$ws.Range("H136").Value = 11820.25
$ws.Range("I136").Value = 2151.25
$ws.Range("J136").Value = 21489.25
$ws.Range("K136").Value = 6453.75
$ws.Range("L136").Value = 64467.75
$ws.Range("M136").Value = -3903.75
$ws.Range("N136").Value = -69567.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 1073.45
$ws.Range("J33").Value = 1754.4166
$ws.Range("L33").Value = 10526.4996
$ws.Range("N33").Value = -11092.4996

$ws.Range("H131").Value = 809.165
$ws.Range("J131").Value = 853.23254
$ws.Range("L131").Value = 2559.69762
$ws.Range("N131").Value = -12639.69762

$ws.Range("H136").Value = 2643.3333
$ws.Range("I136").Value = 1965
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 5895
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -795
$ws.Range("N136").Value = -22200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 44994
$ws.Range("J64").Value = 44994
$ws.Range("L64").Value = 44994
$ws.Range("N64").Value = -45490

$ws.Range("H67").Value = 44994
$ws.Range("J67").Value = 44994
$ws.Range("L67").Value = 44994
$ws.Range("N67").Value = -46710

$ws.Range("H120").Value = 35011
$ws.Range("J120").Value = 35011
$ws.Range("L120").Value = 35011
$ws.Range("N120").Value = -44687

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 34854
$ws.Range("J69").Value = 34854
$ws.Range("L69").Value = 34854
$ws.Range("N69").Value = -36476

$ws.Range("H72").Value = 34854
$ws.Range("J72").Value = 34854
$ws.Range("L72").Value = 104562
$ws.Range("N72").Value = -112674

$ws.Range("H132").Value = 2383.923
$ws.Range("I132").Value = 2401.34
$ws.Range("J132").Value = 2325.8667
$ws.Range("K132").Value = 7204.02
$ws.Range("L132").Value = 6977.6001
$ws.Range("M132").Value = -4674.02
$ws.Range("N132").Value = -12037.6001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 5780.6665
$ws.Range("I41").Value = 5342
$ws.Range("K41").Value = 5342
$ws.Range("M41").Value = -4952

$ws.Range("H69").Value = 13617.429
$ws.Range("J69").Value = 13617.429
$ws.Range("L69").Value = 13617.429
$ws.Range("N69").Value = -15115.429

$ws.Range("H72").Value = 13617.429
$ws.Range("J72").Value = 13617.429
$ws.Range("L72").Value = 40852.287
$ws.Range("N72").Value = -48340.287

$ws.Range("H121").Value = 39999
$ws.Range("J121").Value = 39999
$ws.Range("L121").Value = 39999
$ws.Range("N121").Value = -43493
